$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column numeric-looking text values are written as text (preserve exact literal formatting)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Rows 35/36: Frax and HuobiToken swap places (rank order changes)
$ws.Range("D2").Value = '26.672.50'
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").Value = '1.802.67'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").Value = '308.92'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").Value = '0.9979'
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").Value = '0.4286'
$ws.Range("E7").Value = '  +2.12%  '
$ws.Range("D8").Value = '0.3639'
$ws.Range("E8").Value = '  +1.43%  '
$ws.Range("D9").Value = '0.07191'
$ws.Range("E9").Value = '  +1.19%  '
$ws.Range("D10").Value = '0.8607'
$ws.Range("E10").Value = '  +1.93%  '
$ws.Range("D11").Value = '20.73'
$ws.Range("E11").Value = '  +2.95%  '
$ws.Range("D12").Value = '1.896.12'
$ws.Range("E12").Value = '  +5.22%  '
$ws.Range("D13").Value = '6.580'
$ws.Range("E13").Value = '  +3.42%  '
$ws.Range("D14").Value = '5.306'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = '0.06901'
$ws.Range("E15").Value = '  +2.17%  '
$ws.Range("D16").Value = '0.9989'
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("D17").Value = '79.90'
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("D18").Value = '0.000008852'
$ws.Range("E18").Value = '  +1.77%  '
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = '15.16'
$ws.Range("E20").Value = '  +1.09%  '
$ws.Range("D21").Value = '26.683.65'
$ws.Range("E21").Value = '  -1.39%  '
$ws.Range("D22").Value = '5.155'
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("D23").Value = '11.07'
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("D24").Value = '2.109.76'
$ws.Range("E24").Value = '  +5.13%  '
$ws.Range("D25").Value = '152.02'
$ws.Range("E25").Value = '  -0.49%  '
$ws.Range("D26").Value = '1.838'
$ws.Range("E26").Value = '  -4.40%  '
$ws.Range("D27").Value = '18.21'
$ws.Range("E27").Value = '  +0.70%  '
$ws.Range("D28").Value = '5.180'
$ws.Range("E28").Value = '  +3.41%  '
$ws.Range("D29").Value = '1.892'
$ws.Range("E29").Value = '  +15.18%  '
$ws.Range("D30").Value = '114.95'
$ws.Range("E30").Value = '  +1.79%  '
$ws.Range("D31").Value = '0.08919'
$ws.Range("E31").Value = '  -1.02%  '
$ws.Range("D32").Value = '0.7486'
$ws.Range("E32").Value = '  +3.51%  '
$ws.Range("D33").Value = '1.160'
$ws.Range("E33").Value = '  +6.99%  '
$ws.Range("D34").Value = '4.393'
$ws.Range("E34").Value = '  +1.80%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.752'
$ws.Range("E35").Value = '  -3.76%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = '1.000'
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").Value = '1.119'
$ws.Range("E37").Value = '  +3.74%  '
$ws.Range("D38").Value = '0.05170'
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("D39").Value = '0.01907'
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("D40").Value = '0.5001'
$ws.Range("E40").Value = '  +0.84%  '
$ws.Range("D41").Value = '0.1628'
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").Value = '2.618'
$ws.Range("E42").Value = '  +0.81%  '
$ws.Range("D43").Value = '6.449'
$ws.Range("E43").Value = '  +9.09%  '
$ws.Range("D44").Value = '8.234'
$ws.Range("E44").Value = '  +2.42%  '
$ws.Range("D45").Value = '106.13'
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("D46").Value = '10.25'
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("D47").Value = '0.9984'
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").Value = '1.645'
$ws.Range("E48").Value = '  +2.71%  '
$ws.Range("D49").Value = '0.4545'
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("D50").Value = '0.06232'
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("D51").Value = '1.786'
$ws.Range("E51").Value = '  +4.94%  '
